# Regenerate the "K" column (column G) values for each data row.
# These values are the (re)computed strikeout-style counts ("K" instead of
# the old "Strike#" naming) produced after regenerating std/mean and the
# s_vals calculation pipeline. We simply write the newly computed scalar
# values back into column G for rows 2-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 1
    11 = 3
    12 = 0
    13 = 2
    14 = 1
    15 = 3
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
